# Auto-generated edit script applying the cryptos.xlsx refresh diff
# (GitHub Actions data refresh, Tue Aug 15 02:59:38 UTC 2023).
#
# Updates the Price (D) and Volume(1h) (E) columns for every coin row,
# plus a 3-row reshuffle at the bottom of the table where BabyDogeCoin
# moves up to rank 44 (row 46), pushing Aptos and RenderToken down one
# row each (rows 47-48), together with their refreshed price/volume.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.350.86'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = '1.841.67'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = "'239.95"
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('D6').Value = "'0.6302"
$ws.Range('E6').Value = '  +0.37%  '
$ws.Range('D7').Value = "'1.000"
$ws.Range('D8').Value = "'0.07461"
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('D9').Value = "'0.2897"
$ws.Range('E9').Value = '  +0.15%  '
$ws.Range('D10').Value = "'24.88"
$ws.Range('E10').Value = '  +2.51%  '
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('D12').Value = '1.840.49'
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('D13').Value = "'4.975"
$ws.Range('E13').Value = '  -0.27%  '
$ws.Range('D14').Value = "'0.6770"
$ws.Range('E14').Value = '  +0.06%  '
$ws.Range('D15').Value = "'0.00001031"
$ws.Range('E15').Value = '  +2.11%  '
$ws.Range('D16').Value = "'81.87"
$ws.Range('E16').Value = '  -0.16%  '
$ws.Range('E17').Value = '  +1.39%  '
$ws.Range('D18').Value = '29.326.96'
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('D19').Value = "'229.09"
$ws.Range('E19').Value = '  +0.40%  '
$ws.Range('D20').Value = "'12.32"
$ws.Range('E20').Value = '  +0.42%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('D22').Value = "'7.388"
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').Value = "'158.06"
$ws.Range('E24').Value = '  -0.27%  '
$ws.Range('D25').Value = "'8.525"
$ws.Range('E25').Value = '  +1.42%  '
$ws.Range('D26').Value = "'0.1350"
$ws.Range('E26').Value = '  -1.83%  '
$ws.Range('D27').Value = "'17.45"
$ws.Range('E27').Value = '  -0.56%  '
$ws.Range('D28').Value = "'0.06905"
$ws.Range('E28').Value = '  +7.44%  '
$ws.Range('D29').Value = "'1.457"
$ws.Range('E29').Value = '  +4.88%  '
$ws.Range('D30').Value = "'1.481"
$ws.Range('E30').Value = '  +0.45%  '
$ws.Range('D31').Value = "'4.069"
$ws.Range('E31').Value = '  +0.61%  '
$ws.Range('D32').Value = "'4.058"
$ws.Range('E32').Value = '  -0.47%  '
$ws.Range('D33').Value = "'1.832"
$ws.Range('E33').Value = '  +0.74%  '
$ws.Range('D34').Value = "'1.139"
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').Value = "'0.7007"
$ws.Range('E35').Value = '  +0.87%  '
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('D37').Value = "'0.01845"
$ws.Range('E37').Value = '  +1.68%  '
$ws.Range('D38').Value = "'2.818"
$ws.Range('E38').Value = '  -0.50%  '
$ws.Range('D39').Value = '1.238.39'
$ws.Range('E39').Value = '  -0.43%  '
$ws.Range('D40').Value = "'6.779"
$ws.Range('E40').Value = '  +3.88%  '
$ws.Range('D41').Value = "'0.9415"
$ws.Range('E41').Value = '  +3.47%  '
$ws.Range('D42').Value = "'0.9989"
$ws.Range('E42').Value = '  +0.10%  '
$ws.Range('D43').Value = '1.989.76'
$ws.Range('E43').Value = '  -0.68%  '
$ws.Range('D44').Value = "'101.01"
$ws.Range('E44').Value = '  -0.25%  '
$ws.Range('D45').Value = "'65.35"
$ws.Range('E45').Value = '  -1.19%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = "'0.00000000120"
$ws.Range('E46').Value = '  +3.18%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').Value = "'7.036"
$ws.Range('E47').Value = '  -0.16%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = "'1.709"
$ws.Range('E48').Value = '  +2.70%  '
$ws.Range('D49').Value = "'8.979"
$ws.Range('E49').Value = '  -0.63%  '
$ws.Range('D50').Value = "'0.1144"
$ws.Range('E50').Value = '  -1.83%  '
$ws.Range('D51').Value = "'0.3912"
$ws.Range('E51').Value = '  -0.73%  '
